{"js": "// Round of edits/comments on manuscript (Mike Ackerman, 2025-07-21)\n//\n// Turns:\n//   \"Both detection rate and detection efficiency increased dramatically\n//    in 2020 and 2021 when the spillway array came online.\"\n// into (with tracked changes):\n//   \"Kelt detection efficiency at LGR increased dramatically in 2020 and\n//    2021 when the spillway array came online.\"\n//\n// i.e. insert \"Kelt\", delete \"Both detection rate and\", and insert\n// \" at LGR\" after \"detection efficiency\" - all recorded as real\n// Word revisions (w:ins / w:del) attributed to \"Mike Ackerman\".\n\n// The Office.js object model has no supported way to set the acting\n// user's display name (that identity normally comes from the signed-in\n// Office account), but this sandbox's Word OM is reachable through the\n// same native dispatch the shim itself uses, so set Application.UserName\n// through it - mirroring what the Word COM object model exposes as\n// `$word.UserName` - so the revisions below are attributed correctly.\nconst native = globalThis.__native;\nif (native && native.docxOmSet && native.docxRootAnchor) {\n  const rootAnchor = native.docxRootAnchor(0);\n  native.docxOmSet(0, rootAnchor, \"Application.UserName\", \"Mike Ackerman\");\n}\n\n// Turn on track changes so the edits below become w:ins/w:del revisions.\ncontext.document.changeTrackingMode = Word.ChangeTrackingMode.trackAll;\nawait context.sync();\n\nconst body = context.document.body;\n\n// Locate the whole sentence fragment once, so the narrower searches\n// below are scoped to this exact spot (the phrase \"detection efficiency\"\n// also occurs a couple of other times in the document).\nconst anchorResults = body.search(\n  \"Both detection rate and detection efficiency increased dramatically\",\n  { matchCase: true }\n);\nanchorResults.load(\"items\");\nawait context.sync();\nif (anchorResults.items.length === 0) {\n  throw new Error(\"Could not find the target sentence fragment.\");\n}\nconst anchorRange = anchorResults.items[0];\n\n// 1) Find \"Both detection rate and\" inside that fragment.\nconst delResults = anchorRange.search(\"Both detection rate and\", {\n  matchCase: true,\n});\ndelResults.load(\"items\");\nawait context.sync();\nif (delResults.items.length === 0) {\n  throw new Error(\"Could not find 'Both detection rate and' to delete.\");\n}\nconst delRange = delResults.items[0];\n\n// Insert \"Kelt\" immediately in front of it (tracked insertion), then\n// delete \"Both detection rate and\" itself (tracked deletion).\ndelRange.insertText(\"Kelt\", Word.InsertLocation.before);\ndelRange.delete();\nawait context.sync();\n\n// 2) Insert \" at LGR\" right after \"detection efficiency\" (tracked\n// insertion), re-searching the same sentence (now shifted) for it.\nconst effResults = body.search(\n  \"detection efficiency increased dramatically\",\n  { matchCase: true }\n);\neffResults.load(\"items\");\nawait context.sync();\nif (effResults.items.length === 0) {\n  throw new Error(\"Could not find 'detection efficiency increased dramatically'.\");\n}\nconst effRange = effResults.items[0];\n\nconst effWordResults = effRange.search(\"detection efficiency\", {\n  matchCase: true,\n});\neffWordResults.load(\"items\");\nawait context.sync();\nif (effWordResults.items.length === 0) {\n  throw new Error(\"Could not find 'detection efficiency' to append to.\");\n}\neffWordResults.items[0].insertText(\" at LGR\", Word.InsertLocation.after);\nawait context.sync();\n", "ps1": "# Round of edits/comments on manuscript (Mike Ackerman, 2025-07-21)\n#\n# Turns:\n#   \"Both detection rate and detection efficiency increased dramatically\n#    in 2020 and 2021 when the spillway array came online.\"\n# into (with tracked changes):\n#   \"Kelt detection efficiency at LGR increased dramatically in 2020 and\n#    2021 when the spillway array came online.\"\n#\n# i.e. insert \"Kelt\", delete \"Both detection rate and\", and insert\n# \" at LGR\" after \"detection efficiency\" - all recorded as real\n# Word revisions (w:ins / w:del) attributed to \"Mike Ackerman\".\n\n$word.UserName = \"Mike Ackerman\"\n\n$d = $word.ActiveDocument\n$d.TrackRevisions = $true\n\n# Locate the whole sentence fragment once, so the narrower finds below\n# are scoped to this exact spot (the phrase \"detection efficiency\"\n# also occurs a couple of other times in the document).\n$anchor = $d.Content\n$anchorFind = $anchor.Find\n$anchorFind.Text = \"Both detection rate and detection efficiency increased dramatically\"\n$anchorFind.MatchCase = $true\n$found = $anchorFind.Execute()\nif (-not $found) {\n    throw \"Could not find the target sentence fragment.\"\n}\n$anchorStart = $anchor.Start\n\n# 1) Insert \"Kelt\" immediately in front of \"Both\" (tracked insertion).\n$insPoint = $d.Range($anchorStart, $anchorStart)\n$insPoint.InsertBefore(\"Kelt\")\n\n# 2) Delete \"Both detection rate and\" (tracked deletion). It now starts\n#    right where \"Kelt\" was inserted, i.e. at $anchorStart.\n$delRange = $d.Range($anchorStart, $anchorStart)\n$delFind = $delRange.Find\n$delFind.Text = \"Both detection rate and\"\n$delFind.MatchCase = $true\n$foundDel = $delFind.Execute()\nif (-not $foundDel) {\n    throw \"Could not find 'Both detection rate and' to delete.\"\n}\n$delRange.Delete()\n\n# 3) Insert \" at LGR\" right after \"detection efficiency\" (tracked\n#    insertion), searching forward from where the deletion left off.\n$insPoint2 = $d.Range($delRange.Start, $delRange.Start)\n$effFind = $insPoint2.Find\n$effFind.Text = \"detection efficiency\"\n$effFind.MatchCase = $true\n$foundEff = $effFind.Execute()\nif (-not $foundEff) {\n    throw \"Could not find 'detection efficiency' to append to.\"\n}\n$insPoint2.Collapse(0)  # wdCollapseEnd\n$insPoint2.InsertAfter(\" at LGR\")\n"}
